$wb = $excel.ActiveWorkbook
$runs = $wb.Worksheets.Item("runs")
$params = $wb.Worksheets.Item("params")

# --- "runs" sheet: Stop Run ID (inclusive) 5 -> 15 ---
$runs.Activate()
$runs.Range("B2").Value = 15
$runs.Range("B1").Select()

# --- "params" sheet: update existing rows 2-6 ---
$params.Cells.Item(2,2).Value = 10   # B2: 20 -> 10

$params.Cells.Item(3,2).Value = 10   # B3: 0 -> 10
$params.Cells.Item(3,3).Value = 0    # C3: 20 -> 0
$params.Cells.Item(3,11).Value = "config1"

$params.Cells.Item(4,3).Value = 0    # C4: 10 -> 0
$params.Cells.Item(4,11).Value = "config1"

$params.Cells.Item(5,2).Value = 10   # B5: 15 -> 10
$params.Cells.Item(5,3).Value = 0    # C5: 15 -> 0
$params.Cells.Item(5,11).Value = "config1"

$params.Cells.Item(6,2).Value = 10   # B6: 0 -> 10
$params.Cells.Item(6,3).Value = 0    # C6: 30 -> 0
$params.Cells.Item(6,11).Value = "config1"

# --- new rows 7-16 ---
$newRows = @(
  @(6, 0, 10, 1,   3, 270, 5, 3, 5, 500, "config1"),
  @(7, 0, 10, 0.9, 3, 270, 5, 3, 5, 500, "config1"),
  @(8, 0, 10, 0.8, 3, 270, 5, 3, 5, 500, "config1"),
  @(9, 0, 10, 0.7, 3, 270, 5, 3, 5, 500, "config1"),
  @(10,0, 10, 0.6, 3, 270, 5, 3, 5, 500, "config1"),
  @(11,5, 5,  1,   3, 270, 5, 3, 5, 500, "config1"),
  @(12,5, 5,  0.9, 3, 270, 5, 3, 5, 500, "config1"),
  @(13,5, 5,  0.8, 3, 270, 5, 3, 5, 500, "config1"),
  @(14,5, 5,  0.7, 3, 270, 5, 3, 5, 500, "config1"),
  @(15,5, 5,  0.6, 3, 270, 5, 3, 5, 500, "config1")
)

$r = 7
foreach ($row in $newRows) {
  for ($c = 1; $c -le 11; $c++) {
    $params.Cells.Item($r, $c).Value = $row[$c-1]
  }
  $r++
}
